$d = $word.ActiveDocument

# --- Add a leading "//" comment-marker to the paragraphs that are missing
#     one (several of the time-parsing examples gained the same "//" prefix
#     the other lines in this scratch file already used). ---
$prefixParagraphs = @(1, 2, 8, 9, 10, 11, 12)
foreach ($idx in $prefixParagraphs) {
    $d.Paragraphs.Item($idx).Range.InsertBefore("//")
}

# --- The last paragraph loses its "//" marker (it is now prose, not a
#     commented-out example), while the "_GoBack" bookmark - which used to
#     sit at the end of the "дцадцать 1 двадцать 1 - 20:01" paragraph - is
#     now the most-recently-edited spot, at the very start of the last
#     paragraph. Move it there, then strip the leading "//". ---
$last = $d.Paragraphs.Last
$lastStart = $last.Range.Duplicate
$lastStart.Collapse(1)

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$d.Bookmarks.Add("_GoBack", $lastStart)

$last.Range.Find.Execute("//", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
